$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" column (column G) values, replacing old Strike# values with
# regenerated K values (regen save_data to use K instead of Strike#).
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 2
    12 = 3
    13 = 0
    14 = 1
    15 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 0
    22 = 2
    23 = 0
    24 = 1
    25 = 2
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 2
    33 = 0
    34 = 3
    35 = 2
    36 = 2
    37 = 0
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 3
    43 = 1
    44 = 0
    45 = 0
    46 = 0
    47 = 1
    48 = 2
    49 = 0
    50 = 1
    51 = 2
    52 = 2
    53 = 0
    54 = 3
    55 = 3
    56 = 1
    57 = 3
    58 = 0
    59 = 0
    60 = 1
    61 = 2
    62 = 2
    63 = 2
    64 = 2
    66 = 2
    67 = 3
    68 = 1
    69 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
